$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string text updates -------------------------------------------------

# "• Never / • Rarely / • Some of the time / • Most of the time" -> append "  • NA"
# (shared by B2 and B3)
$answers1 = "• Never `n• Rarely `n• Some of the time`n• Most of the time                               • NA"
$ws.Range("B2").Value = $answers1
$ws.Range("B3").Value = $answers1

# "15, 34" -> "15, 34, 38" (shared by E2, E4, E5, E6, E7)
$ws.Range("E2").Value = "15, 34, 38"
$ws.Range("E4").Value = "15, 34, 38"
$ws.Range("E5").Value = "15, 34, 38"
$ws.Range("E6").Value = "15, 34, 38"
$ws.Range("E7").Value = "15, 34, 38"

# "15. 34" -> "15. 34, 38" (E3 only)
$ws.Range("E3").Value = "15. 34, 38"

# "•  Yes / •  No" -> append "  • NA" (shared by B4, B5, B6, B7)
$answers2 = "•  Yes`n•  No                                                       • NA"
$ws.Range("B4").Value = $answers2
$ws.Range("B5").Value = $answers2
$ws.Range("B6").Value = $answers2
$ws.Range("B7").Value = $answers2

# --- Sheet view / selection updates ---------------------------------------------
# Remove the scrolled-down top-left cell (was A5) and move the selection from
# E7 to B7.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("B7").Select() | Out-Null
